$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the active cell/selection that appears in the saved file
$ws.Range("E8").Select()
